$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 634, shifting the existing rows 634-672
# (and everything below) down by two rows.
$ws.Range("A634:A635").EntireRow.Insert()

# Row 634: new "Primera" quality record for the week of 2022-01-24 (44585)
$ws.Range("A634").Value = 6
$ws.Range("B634").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C634").Value = "Metropolitana"
$ws.Range("D634").Value = 44585
$ws.Range("E634").Value = 13
$ws.Range("F634").Value = 100112009
$ws.Range("G634").Value = "Acelga"
$ws.Range("H634").Value = "Sin especificar"
$ws.Range("I634").Value = "Primera"
$ws.Range("J634").Value = 100
$ws.Range("K634").Value = 16000
$ws.Range("L634").Value = 16000
$ws.Range("M634").Value = 16000
$ws.Range("N634").Value = "`$/docena de atados"
$ws.Range("O634").Value = "Región Metropolitana"
$ws.Range("P634").Value = 5333
$ws.Range("Q634").Value = 3
$ws.Range("R634").Value = "Hortaliza"

# Row 635: new "Segunda" quality record for the same week (44585)
$ws.Range("A635").Value = 6
$ws.Range("B635").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C635").Value = "Metropolitana"
$ws.Range("D635").Value = 44585
$ws.Range("E635").Value = 13
$ws.Range("F635").Value = 100112009
$ws.Range("G635").Value = "Acelga"
$ws.Range("H635").Value = "Sin especificar"
$ws.Range("I635").Value = "Segunda"
$ws.Range("J635").Value = 70
$ws.Range("K635").Value = 12000
$ws.Range("L635").Value = 12000
$ws.Range("M635").Value = 12000
$ws.Range("N635").Value = "`$/docena de atados"
$ws.Range("O635").Value = "Región Metropolitana"
$ws.Range("P635").Value = 4000
$ws.Range("Q635").Value = 3
$ws.Range("R635").Value = "Hortaliza"

# Make sure the new date cells use the same date/time number format as the
# rest of column D.
$ws.Range("D634:D635").NumberFormat = $ws.Range("D636").NumberFormat()
